$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Meetups")

# Insert a new row before row 7; this shifts the existing rows 7-17 down to 8-18
# and Excel copies row 6's formatting into the new row 7 (matches the target
# diff: A7 keeps style 1, B7/C7 keep style 3).
$ws.Rows.Item(7).Insert()

# New row 7: "R Package Development" talk (2022-03-01)
$ws.Range("A7").Value = 44621
$ws.Range("B7").Value = "7:00 pm"
$ws.Range("C7").Value = "8:30 pm"
$ws.Range("D7").Value = "R Package Development"
$ws.Range("F7").Value = "vlyVKGSVCsk"
$ws.Range("G7").Value = "2022-03-01-R_Package_Development"

# H7 gets a fresh (non-inherited) style before wrapping so it lands on the
# default font (matches the new cellXfs entry: fontId=0 + wrapText=1).
$ws.Range("H7").Style = "Normal"
$ws.Range("H7").WrapText = $true
$ws.Range("H7").Value = "R script located here: https://github.com/jbryer/DATA606Spring2022/blob/main/Slides/2022-03-01-R_Package_Development/Build_R_Package.R" + [char]10 + "You can download the supporting materials here: https://github.com/jbryer/DATA606Spring2022/blob/main/Slides/2022-03-01-R_Package_Development/"

# Row height for the new row (wrapped text needs extra height)
$ws.Rows.Item(7).RowHeight = 23

# Column G needs to be wider to accommodate the new slide-deck folder name
$ws.Columns.Item(7).ColumnWidth = 32.83

# Update the active selection to match where the author ended up
$ws.Range("G9").Select()
